# The edit re-orders the 5 observation records currently sitting in rows
# 2-6 of the "Artfynd" sheet: the record in row 2 (Id 80923932, "Ask")
# moves down to row 6, and every other record shifts up by one row
# (row3->row2, row4->row3, row5->row4, row6->row5). All other sheet
# content (header row, column layout, styles) is untouched.
#
# Text values are assigned with a leading "'" (apostrophe) so the engine
# never reinterprets look-alike dates/numbers (e.g. "2014-07-24",
# "00:00") as real date/time serials, then the style is reset back to
# "Normal" so no stray number-format / quote-prefix formatting is left
# behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($a1, $text) {
    $ws.Range($a1).Value = "'" + $text
    $ws.Range($a1).Style = "Normal"
}

function Set-EmptyText($a1) {
    $ws.Range($a1).Value = "'"
    $ws.Range($a1).Style = "Normal"
}

# ---- Row 2  (becomes the old row-3 record: Skogsnycklar) ----
$ws.Range("A2").Value = 81368550
$ws.Range("B2").Value = 96252
Set-Text "D2" "LC"
$ws.Range("E2").Value = 223591
Set-Text "F2" "Skogsnycklar"
Set-Text "G2" "Dactylorhiza maculata subsp. fuchsii"
Set-Text "H2" "(Druce) Hyl."
Set-Text "P2" "1 km SO om Södergård, Sk"
$ws.Range("Q2").Value = 420859.7601243296
$ws.Range("R2").Value = 6190239.832464891
$ws.Range("S2").Value = 50
Set-Text "Y2" "2014-07-24"
Set-Text "AA2" "2014-07-24"
Set-Text "AI2" "Fuktigt hässle"
Set-EmptyText "AR2"
Set-Text "AW2" "Charlotte Wigermo"
Set-Text "AX2" "Torbjörn Tyler"
Set-Text "AY2" "Skånes Flora Millora 2008-2015"

# ---- Row 3  (becomes the old row-4 record: Grönvit nattviol) ----
$ws.Range("A3").Value = 81368358
$ws.Range("B3").Value = 96370
$ws.Range("E3").Value = 219875
Set-Text "F3" "Grönvit nattviol"
Set-Text "G3" "Platanthera chlorantha"
Set-Text "H3" "(Custer) Rchb."
$ws.Range("AR3").ClearContents()

# ---- Row 4  (becomes the old row-5 record: Blåsippa) ----
$ws.Range("A4").Value = 81368477
$ws.Range("B4").Value = 98520
$ws.Range("E4").Value = 222498
Set-Text "F4" "Blåsippa"
Set-Text "G4" "Hepatica nobilis"
Set-Text "H4" "Schreb."

# ---- Row 5  (becomes the old row-6 record: Skogsknipprot) ----
$ws.Range("A5").Value = 81368525
$ws.Range("B5").Value = 96312
$ws.Range("E5").Value = 219798
Set-Text "F5" "Skogsknipprot"
Set-Text "G5" "Epipactis helleborine"
Set-Text "H5" "(L.) Crantz"
$ws.Range("Q5").Value = 420758.613196145
$ws.Range("R5").Value = 6190338.232892067
Set-Text "AI5" "Rik ek-/hasselskog"

# ---- Row 6  (becomes the old row-2 record: Ask) ----
$ws.Range("A6").Value = 80923932
$ws.Range("B6").Value = 103813
Set-Text "D6" "EN"
$ws.Range("E6").Value = 220785
Set-Text "F6" "Ask"
Set-Text "G6" "Fraxinus excelsior"
Set-Text "H6" "L."
Set-Text "P6" "Tågarp, Sk"
$ws.Range("Q6").Value = 420921.7954929644
$ws.Range("R6").Value = 6190309.99155066
$ws.Range("S6").Value = 10
Set-Text "Y6" "2019-07-29"
Set-Text "AA6" "2019-07-29"
Set-Text "AI6" "Lövskog"
Set-Text "AW6" "Örjan Fritz"
Set-Text "AX6" "Örjan Fritz"
Set-EmptyText "AY6"
